# Update the cryptocurrency price/volume table with the latest scraped
# values. Columns: D = Price, E = Volume(1h). Values are stored as text,
# so number formats are forced to "@" (text) before assignment to avoid
# Excel auto-converting numeric-looking strings (e.g. "215.45") into
# real numbers, which would change both their type and their formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    # Reset back to the default (un-styled) cell style so we don't leave
    # an extraneous explicit number format on cells that originally had
    # none, keeping the sheet's formatting identical to before the edit.
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "25.756.34"
Set-TextValue $ws.Range("E2") "  -0.19%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.633.86"
Set-TextValue $ws.Range("E3") "  -0.11%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  -0.11%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "215.45"
Set-TextValue $ws.Range("E5") "  +0.09%  "

# Row 6 - XRP
Set-TextValue $ws.Range("E6") "  -0.61%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("E8") "  -0.70%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("E9") "  -1.37%  "

# Row 10 - Solana
Set-TextValue $ws.Range("D10") "19.54"
Set-TextValue $ws.Range("E10") "  -1.61%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.0792"
Set-TextValue $ws.Range("E11") "  +0.85%  "

# Row 12 - Polkadot
Set-TextValue $ws.Range("E12") "  +0.48%  "

# Row 13 - Wrapped liquid staked Ether 2.0
Set-TextValue $ws.Range("D13") "1.859.41"
Set-TextValue $ws.Range("E13") "  -0.10%  "

# Row 14 - Wrapped Ether
Set-TextValue $ws.Range("D14") "1.635.12"
Set-TextValue $ws.Range("E14") "  -0.03%  "

# Row 15 - Polygon
Set-TextValue $ws.Range("E15") "  +0.22%  "

# Row 16 - Shiba Inu
Set-TextValue $ws.Range("E16") "  -1.11%  "

# Row 17 - Litecoin
Set-TextValue $ws.Range("D17") "63.14"
Set-TextValue $ws.Range("E17") "  +0.11%  "

# Row 18 - Wrapped BTC
Set-TextValue $ws.Range("D18") "25.777.22"
Set-TextValue $ws.Range("E18") "  -0.15%  "

# Row 19 - Dai
Set-TextValue $ws.Range("E19") "  -0.08%  "

# Row 20 - Uniswap
Set-TextValue $ws.Range("D20") "4.45"
Set-TextValue $ws.Range("E20") "  +0.23%  "

# Row 21 - Bitcoin Cash
Set-TextValue $ws.Range("D21") "192.50"
Set-TextValue $ws.Range("E21") "  -0.79%  "

# Row 22 - Avalanche
Set-TextValue $ws.Range("E22") "  +0.23%  "

# Row 23 - Chainlink
Set-TextValue $ws.Range("E23") "  +1.96%  "

# Row 24 - BinanceUSD
Set-TextValue $ws.Range("E24") "  -0.04%  "

# Row 25 - Toncoin
Set-TextValue $ws.Range("E25") "  +1.99%  "

# Row 26 - Monero
Set-TextValue $ws.Range("D26") "142.69"
Set-TextValue $ws.Range("E26") "  +2.53%  "

# Row 27 - Stellar
Set-TextValue $ws.Range("E27") "  +1.55%  "

# Row 28 - Cosmos
Set-TextValue $ws.Range("E28") "  +0.57%  "

# Row 30 - PancakeSwap
Set-TextValue $ws.Range("E30") "  +0.00%  "

# Row 31 - Hedera
Set-TextValue $ws.Range("E31") "  -1.17%  "

# Row 32 - Internet Computer (DFINITY)
Set-TextValue $ws.Range("D32") "3.34"
Set-TextValue $ws.Range("E32") "  +0.03%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("E33") "  -1.05%  "

# Row 34 - Lido DAO Token
Set-TextValue $ws.Range("E34") "  -1.36%  "

# Row 35 - Huobi Token
Set-TextValue $ws.Range("E35") "  -0.21%  "

# Row 36 - ARBITRUM
Set-TextValue $ws.Range("E36") "  +0.82%  "

# Row 37 - Maker
Set-TextValue $ws.Range("D37") "1.132.23"
Set-TextValue $ws.Range("E37") "  +2.34%  "

# Row 38 - MX Token
Set-TextValue $ws.Range("E38") "  -2.19%  "

# Row 39 - Immutable X
Set-TextValue $ws.Range("E39") "  -1.08%  "

# Row 40 - VeChain
Set-TextValue $ws.Range("E40") "  -0.85%  "

# Row 41 - Pax Dollar
Set-TextValue $ws.Range("E41") "  +0.09%  "

# Row 42 - mCoin
Set-TextValue $ws.Range("E42") "  +0.25%  "

# Row 43 - Frax Share
Set-TextValue $ws.Range("D43") "5.57"
Set-TextValue $ws.Range("E43") "  -0.23%  "

# Row 44 - Quant
Set-TextValue $ws.Range("D44") "100.74"

# Row 45 - Trust Wallet Token
Set-TextValue $ws.Range("E45") "  -0.42%  "

# Row 46 - Rocket Pool ETH
Set-TextValue $ws.Range("D46") "1.768.47"
Set-TextValue $ws.Range("E46") "  +0.03%  "

# Row 47 - Baby Doge Coin
Set-TextValue $ws.Range("E47") "  +3.98%  "

# Row 48 - Aave
Set-TextValue $ws.Range("D48") "55.37"
Set-TextValue $ws.Range("E48") "  -0.41%  "

# Row 49 - Cronos
Set-TextValue $ws.Range("D49") "0.0508"
Set-TextValue $ws.Range("E49") "  +0.98%  "

# Row 50 - Mantle
Set-TextValue $ws.Range("E50") "  -0.34%  "

# Row 51 - Render Token
Set-TextValue $ws.Range("E51") "  +3.22%  "
